# "fix bug of chinese name"
# Adds two new configuration columns (E, F) to the "配置" sheet that hold
# per-class regex helpers used to clean up mis-typed Chinese class names:
#   E = "错误名字去除字符列表"  (characters/patterns to strip from a wrong name)
#   F = "错误名字指定转化"      (explicit replacement mapping for a wrong name)
# Rows 2-5 correspond to the same classes already listed in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "配置" (configuration) sheet

# --- New header row (row 1) ---
$ws.Range("E1").Value = "错误名字去除字符列表"
$ws.Range("F1").Value = "错误名字指定转化"

# --- New per-class values (column E only; F stays empty like the diff) ---
# Row 2/3 => "20计网" classes: no special cleanup pattern, just the class name.
$ws.Range("E2").Value = "20计网"
$ws.Range("E3").Value = "20计网"

# Row 5 => "22机电4" class: regex used to strip stray characters from names.
# (populated before row 4 so new shared-string entries land in the same
#  order as the authoritative edit)
$ws.Range("E5").Value = "(（?机电一体化[4四]?班～?）?)|(22级)|(机电22)|(机电[4四]?班?)"

# Row 4 => "21电商2" class: regex used to strip stray characters from names.
$ws.Range("E4").Value = "21电商[二2]?班?,"

# --- Column sizing for the two new columns ---
$ws.Columns.Item(5).ColumnWidth = 39.45
$ws.Columns.Item(6).ColumnWidth = 37.04

# --- View state: zoom in, and leave the selection on the last edited cell ---
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("E5").Select()
